$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Net loss allocation values for columns F and G, rows 4-13
$data = @{
    4  = @(3751520.25, 2108149.75)
    5  = @(3751520.25, 2108149.75)
    6  = @(2248482.25, 1263524.3799999999)
    7  = @(653114.38, 367014.66)
    8  = @(39164936, 22008558)
    9  = @(14765301, 8297294)
    10 = @(299133.31, 168096.59)
    11 = @(7803477.5, 4385129)
    12 = @(41202360, 23153478)
    13 = @(16480944, 9261392)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 6).Value = $vals[0]
    $ws.Cells.Item($r, 7).Value = $vals[1]
}

# Column K (11) gets blank cells formatted with the builtin "Percent" style
# (number format 0.00%) for rows 4-13. Build the format once on a scratch
# cell and paste it onto the whole range so only a single new style entry
# is added (instead of one per cell).
$tmpl = $ws.Cells.Item(1, 20)
$tmpl.Style = "Percent"
$tmpl.NumberFormat = "0.00%"

$tmpl.Copy() | Out-Null
$percentRange = $ws.Range("K4:K13")
$percentRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$tmpl.Clear() | Out-Null

# bestFit custom column widths for columns G (7) and I (9) -- stored width 11
$bestFitWidth = 11 - (5 / 6)
$ws.Columns.Item(7).ColumnWidth = $bestFitWidth
$ws.Columns.Item(9).ColumnWidth = $bestFitWidth
